# Generate Report for Handoff
# Updates the handoff batch identifier (GUID), the content hash embedded in
# the xlf file names, and the handoff/target timestamps across all three
# worksheets (Overview, zh-cn, de-de) as well as the matching hyperlink
# display text.

$wb = $excel.ActiveWorkbook

$oldGuid = "ff89f51c-597e-45dd-a1f2-052b7b44a0b9"
$newGuid = "9c8c3139-075d-4ce6-8cbe-a5a95d7b5459"

$oldHash = "5f98d28e64e56112ac38dff50e24e5fe8599043c"
$newHash = "ed16e4810bb4842d64c740e1ee5d32744a40e84b"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("D2").Value = "2016-43-14 06:43:45"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldGuid.md") {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-14 06:43:41"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldGuid.md") {
        $hl.TextToDisplay = "$newGuid.md"
    } elseif ($hl.TextToDisplay -eq "$oldGuid.$oldHash.zh-cn.xlf") {
        $hl.TextToDisplay = "$newGuid.$newHash.zh-cn.xlf"
    }
}

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-14 06:43:45"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldGuid.md") {
        $hl.TextToDisplay = "$newGuid.md"
    } elseif ($hl.TextToDisplay -eq "$oldGuid.$oldHash.de-de.xlf") {
        $hl.TextToDisplay = "$newGuid.$newHash.de-de.xlf"
    }
}
